# Update "想去人数" (F column) figures across sheets, matching output
# generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3172
$ws1.Range("F5").Value = 6863
$ws1.Range("F6").Value = 1971
$ws1.Range("F13").Value = 146
$ws1.Range("F14").Value = 175
$ws1.Range("F15").Value = 33

# 演出 (Performances) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3172
$ws4.Range("F3").Value = 14
$ws4.Range("F6").Value = 6863
$ws4.Range("F7").Value = 1971
$ws4.Range("F14").Value = 146
$ws4.Range("F15").Value = 175
$ws4.Range("F16").Value = 33

$wb.Save()
